$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update question 2 header text (B17)
$ws.Range("B17").Value = "問題２　（　　　）に入れるのに最もよいものを、１・２・３・４から一つ選びなさい。"

# New question data for the '白' question type, repeated across rows 18-23
$qType  = "白"
$qText  = "私は主張は単なる（　　）ではなく、確たる証拠に基づいている。"
$choice1 = "爆発"
$choice2 = "視線"
$choice3 = "推測"
$choice4 = "推移"

18..23 | ForEach-Object {
    $r = $_
    $ws.Range("A$r").Value = $qType
    $ws.Range("B$r").Value = $qText
    $ws.Range("C$r").Value = $choice1
    $ws.Range("D$r").Value = $choice2
    $ws.Range("E$r").Value = $choice3
    $ws.Range("F$r").Value = $choice4
}

# Update active selection to E13
$ws.Range("E13").Select()
